$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 30, shifting existing row 30 (and everything below) down by one.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new data point.
$ws.Cells.Item(30, 1).Value = 4
$ws.Cells.Item(30, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(30, 3).Value = "Los Lagos"
$ws.Cells.Item(30, 4).Value = 44600
$ws.Cells.Item(30, 5).Value = 10
$ws.Cells.Item(30, 6).Value = 100112009
$ws.Cells.Item(30, 7).Value = "Acelga"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 80
$ws.Cells.Item(30, 11).Value = 10000
$ws.Cells.Item(30, 12).Value = 10000
$ws.Cells.Item(30, 13).Value = 10000
$ws.Cells.Item(30, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(30, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(30, 16).Value = 833
$ws.Cells.Item(30, 17).Value = 12
$ws.Cells.Item(30, 18).Value = "Hortaliza"
